$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the row for "SETENAY AKSU" (row 16) - subsequent rows shift up by one.
$ws.Rows.Item(16).Delete()

# Remove the row for "NAİLE DİKEÇ" (now row 21 after the previous delete) -
# subsequent rows shift up by one again.
$ws.Rows.Item(21).Delete()

# Append a new last entry, "VOLKAN İZCİ", after "MEHTAP AKDOĞAN" (now row 21).
$ws.Range("A22").Value = "VOLKAN İZCİ"

# Match the final selection left active on the sheet.
$ws.Range("A26").Select()
